$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.727026462554932
$ws.Range("B1").Value = 2.558655500411987
$ws.Range("C1").Value = 2.999873876571655
$ws.Range("D1").Value = 2.598408937454224
$ws.Range("E1").Value = 0.5397922396659851
